# Apply the updated crypto price/volume figures to Sheet1 (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.706.14"
$ws.Range("E2").Value = "  -4.12%  "
$ws.Range("D3").Value = "3.155.71"
$ws.Range("E3").Value = "  -4.33%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'524.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.13%  "
$ws.Range("D6").Value = "'132.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.38%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.156.58"
$ws.Range("E8").Value = "  -4.26%  "
$ws.Range("D9").Value = "'0.453"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.12%  "
$ws.Range("D10").Value = "'7.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.06%  "
$ws.Range("E11").Value = "  -6.83%  "
$ws.Range("D12").Value = "'0.390"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.90%  "
$ws.Range("D13").Value = "3.699.46"
$ws.Range("E13").Value = "  -4.39%  "
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").Value = "'25.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.88%  "
$ws.Range("D16").Value = "3.159.36"
$ws.Range("D17").Value = "57.747.06"
$ws.Range("E17").Value = "  -4.48%  "
$ws.Range("E18").Value = "  -7.66%  "
$ws.Range("D19").Value = "'5.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.72%  "
$ws.Range("D20").Value = "'13.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.00%  "
$ws.Range("D21").Value = "'8.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.76%  "
$ws.Range("D22").Value = "'346.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.09%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'69.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.85%  "
$ws.Range("E25").Value = "  -6.64%  "
$ws.Range("D26").Value = "3.286.31"
$ws.Range("E26").Value = "  -5.09%  "
$ws.Range("D27").Value = "0.0₃0961"
$ws.Range("E27").Value = "  -8.10%  "
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'6.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.64%  "
$ws.Range("D32").Value = "'1.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.10%  "
$ws.Range("E33").Value = "  -9.17%  "
$ws.Range("D34").Value = "'21.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.71%  "
$ws.Range("E35").Value = "  -4.87%  "
$ws.Range("D36").Value = "'4.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.40%  "
$ws.Range("D37").Value = "'159.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.32%  "
$ws.Range("E38").Value = "  -7.50%  "
$ws.Range("E39").Value = "  -7.62%  "
$ws.Range("D40").Value = "'25.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.67%  "
$ws.Range("E41").Value = "  -5.27%  "
$ws.Range("D42").Value = "3.185.58"
$ws.Range("E42").Value = "  -4.39%  "
$ws.Range("D43").Value = "'40.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").Value = "'0.698"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.00%  "
$ws.Range("E45").Value = "  -3.88%  "
$ws.Range("D46").Value = "'3.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.20%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  -8.13%  "
$ws.Range("D49").Value = "2.268.34"
$ws.Range("E49").Value = "  -4.19%  "
$ws.Range("D51").Value = "'20.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.64%  "
